$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 22000
$ws.Range("I21").Value = 50000
$ws.Range("J21").Value = 15000
$ws.Range("K21").Value = 50000
$ws.Range("L21").Value = 15000
$ws.Range("M21").Value = -49532
$ws.Range("N21").Value = -15936

$ws.Range("H23").Value = 22000
$ws.Range("I23").Value = 50000
$ws.Range("J23").Value = 15000
$ws.Range("K23").Value = 50000
$ws.Range("L23").Value = 15000
$ws.Range("M23").Value = -49766
$ws.Range("N23").Value = -15468

$ws.Range("H34").Value = 674959
$ws.Range("I34").Value = 839107.25
$ws.Range("J34").Value = 18366
$ws.Range("K34").Value = 839107.25
$ws.Range("L34").Value = 18366
$ws.Range("M34").Value = -838904.25
$ws.Range("N34").Value = -18772

$ws.Range("H36").Value = 674959
$ws.Range("I36").Value = 839107.25
$ws.Range("J36").Value = 18366
$ws.Range("K36").Value = 839107.25
$ws.Range("L36").Value = 18366
$ws.Range("M36").Value = -838392.25
$ws.Range("N36").Value = -19796

$ws.Range("H75").Value = 312771.34
$ws.Range("J75").Value = 269157
$ws.Range("L75").Value = 269157
$ws.Range("N75").Value = -271029

$ws.Range("H78").Value = 312771.34
$ws.Range("J78").Value = 269157
$ws.Range("L78").Value = 807471
$ws.Range("N78").Value = -816831

$ws.Range("H132").Value = 43395.52
$ws.Range("I132").Value = 43395.52
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 130186.56
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -127656.56
$ws.Range("N132").ClearContents()

$ws.Range("H138").Value = 1842.8182
$ws.Range("I138").Value = 864.4524
$ws.Range("J138").Value = 2563.7192
$ws.Range("K138").Value = 2593.3572
$ws.Range("L138").Value = 7691.1576
$ws.Range("M138").Value = 2546.6428
$ws.Range("N138").Value = -17971.1576


# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20540.715
$ws.Range("I32").Value = 5136.9546
$ws.Range("J32").Value = 112963.27
$ws.Range("K32").Value = 5136.9546
$ws.Range("L32").Value = 112963.27
$ws.Range("M32").Value = -4849.9546
$ws.Range("N32").Value = -113537.27

$ws.Range("H60").Value = 11272.728
$ws.Range("I60").Value = 50000
$ws.Range("J60").Value = 7400
$ws.Range("K60").Value = 50000
$ws.Range("L60").Value = 7400
$ws.Range("M60").Value = -49267
$ws.Range("N60").Value = -8866

$ws.Range("H97").Value = 100000
$ws.Range("I97").Value = 100000
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 100000
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -99504
$ws.Range("N97").ClearContents()

$ws.Range("H110").Value = 58824570
$ws.Range("I110").Value = 83334380
$ws.Range("J110").Value = 1001.2
$ws.Range("K110").Value = 83334380
$ws.Range("L110").Value = 1001.2
$ws.Range("M110").Value = -83332335
$ws.Range("N110").Value = -5091.2

$ws.Range("H122").Value = 1656.2778
$ws.Range("I122").Value = 1351.25
$ws.Range("K122").Value = 4053.75
$ws.Range("M122").Value = -1603.75


# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2249.75
$ws.Range("I20").Value = 1899.6
$ws.Range("J20").Value = 2833.3333
$ws.Range("K20").Value = 1899.6
$ws.Range("L20").Value = 2833.3333
$ws.Range("M20").Value = -1652.6
$ws.Range("N20").Value = -3327.3333

$ws.Range("H99").Value = 1866.35
$ws.Range("I99").Value = 1794.8572
$ws.Range("J99").Value = 2033.1666
$ws.Range("K99").Value = 1794.8572
$ws.Range("L99").Value = 2033.1666
$ws.Range("M99").Value = -296.8571999999999
$ws.Range("N99").Value = -5029.1666

$ws.Range("H134").Value = 4321.619
$ws.Range("I134").Value = 2909.3333
$ws.Range("J134").Value = 6204.6665
$ws.Range("K134").Value = 8727.999899999999
$ws.Range("L134").Value = 18613.9995
$ws.Range("M134").Value = -6192.999899999999
$ws.Range("N134").Value = -23683.9995


# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1420.7391
$ws.Range("I16").Value = 1272.5834
$ws.Range("J16").Value = 1582.3636
$ws.Range("K16").Value = 1272.5834
$ws.Range("L16").Value = 1582.3636
$ws.Range("M16").Value = -985.5834
$ws.Range("N16").Value = -2156.3636

$ws.Range("H45").Value = 5933.4
$ws.Range("I45").Value = 67
$ws.Range("J45").Value = 7400
$ws.Range("K45").Value = 67
$ws.Range("L45").Value = 7400
$ws.Range("M45").Value = 526
$ws.Range("N45").Value = -8586

$ws.Range("H107").Value = 491.93332
$ws.Range("I107").Value = 372.8889
$ws.Range("J107").Value = 670.5
$ws.Range("K107").Value = 372.8889
$ws.Range("L107").Value = 670.5
$ws.Range("M107").Value = 1547.1111
$ws.Range("N107").Value = -4510.5

$ws.Range("H113").Value = 1420.7391
$ws.Range("I113").Value = 1272.5834
$ws.Range("J113").Value = 1582.3636
$ws.Range("K113").Value = 1272.5834
$ws.Range("L113").Value = 1582.3636
$ws.Range("M113").Value = 897.4166
$ws.Range("N113").Value = -5922.3636


# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 8265208.5
$ws.Range("I113").Value = 500
$ws.Range("J113").Value = 9091680
$ws.Range("K113").Value = 1500
$ws.Range("L113").Value = 27275040
$ws.Range("M113").Value = 670
$ws.Range("N113").Value = -27279380

$ws.Range("H131").Value = 6062096
$ws.Range("J131").Value = 7409104
$ws.Range("L131").Value = 22227312
$ws.Range("N131").Value = -22237392


# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4880.1665
$ws.Range("I70").Value = 4955.5815
$ws.Range("J70").Value = 4585.364
$ws.Range("K70").Value = 4955.5815
$ws.Range("L70").Value = 4585.364
$ws.Range("M70").Value = -4685.5815
$ws.Range("N70").Value = -5125.364

$ws.Range("H73").Value = 4880.1665
$ws.Range("I73").Value = 4955.5815
$ws.Range("J73").Value = 4585.364
$ws.Range("K73").Value = 4955.5815
$ws.Range("L73").Value = 4585.364
$ws.Range("M73").Value = -4019.5815
$ws.Range("N73").Value = -6457.364

$ws.Range("H80").Value = 40002864
$ws.Range("I80").Value = 2672.8
$ws.Range("K80").Value = 2672.8
$ws.Range("M80").Value = -1674.8

$ws.Range("H83").Value = 40002864
$ws.Range("I83").Value = 2672.8
$ws.Range("K83").Value = 13364
$ws.Range("M83").Value = -8372

$ws.Range("H113").Value = 2077.75
$ws.Range("I113").Value = 1970.3334
$ws.Range("J113").Value = 2400
$ws.Range("K113").Value = 1970.3334
$ws.Range("L113").Value = 2400
$ws.Range("M113").Value = 199.6666
$ws.Range("N113").Value = -6740


# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 10000344
$ws.Range("I16").Value = 11111460
$ws.Range("K16").Value = 11111460
$ws.Range("M16").Value = -11111290

$ws.Range("H55").Value = 245.65218
$ws.Range("I55").Value = 219.26315
$ws.Range("J55").Value = 371
$ws.Range("K55").Value = 219.26315
$ws.Range("L55").Value = 371
$ws.Range("M55").Value = -46.26315
$ws.Range("N55").Value = -717

$ws.Range("H61").Value = 1132.4166
$ws.Range("I61").Value = 626.2857
$ws.Range("J61").Value = 1841
$ws.Range("K61").Value = 626.2857
$ws.Range("L61").Value = 1841
$ws.Range("M61").Value = -424.2857
$ws.Range("N61").Value = -2245

$ws.Range("H98").Value = 30000
$ws.Range("J98").Value = 30000
$ws.Range("L98").Value = 30000
$ws.Range("N98").Value = -35990

$ws.Range("H113").Value = 1132.4166
$ws.Range("I113").Value = 626.2857
$ws.Range("J113").Value = 1841
$ws.Range("K113").Value = 626.2857
$ws.Range("L113").Value = 1841
$ws.Range("M113").Value = 1543.7143
$ws.Range("N113").Value = -6181


# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 4888.8887
$ws.Range("I54").Value = 4250
$ws.Range("J54").Value = 10000
$ws.Range("K54").Value = 4250
$ws.Range("L54").Value = 10000
$ws.Range("M54").Value = -3730
$ws.Range("N54").Value = -11040

$ws.Range("H104").Value = 25600
$ws.Range("J104").Value = 25600
$ws.Range("L104").Value = 25600
$ws.Range("N104").Value = -32588

$ws.Range("H113").Value = 716.7059
$ws.Range("I113").Value = 875.6923
$ws.Range("J113").Value = 200
$ws.Range("K113").Value = 2627.0769
$ws.Range("L113").Value = 600
$ws.Range("M113").Value = -457.0769
$ws.Range("N113").Value = -4940

